# Traceability matrix update:
#   - a new "Getroffene Designentscheidungen" entry is inserted as a data
#     row right after the "Source Code" section entry (old row 26),
#     pushing the "Source Code" section header and everything below it
#     down by one row.
#   - selection is left on the newly relevant area (F28, which is the
#     "Source Code" section header row after the shift).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row before the old row 27 ("Source Code" section header).
$ws.Rows.Item(27).Insert()

# Populate the new row with the "design decisions" entry, matching the
# layout used by the other data rows (B: short name, C: location,
# D: document name, F: remark/description).
$ws.Range("B27").Value = "Getroffene Designentscheidungen"
$ws.Range("C27").Value = "Mobile/Designentscheidungen"
$ws.Range("D27").Value = "Designentscheidungen.pdf"
$ws.Range("F27").Value = "Getroffene Designentscheidungen im Projekt"

# Reflect the saved selection state from the edited workbook.
$ws.Range("F28").Select()
